$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Table on slide 5: switch the table to a different (built-in) table style.
#    OOXML: <a:tableStyleId>{E69BC347-...}</a:tableStyleId>
#             -> <a:tableStyleId>{8F568475-D16F-4FD6-A70D-62D4BDE3F239}</a:tableStyleId>
# ---------------------------------------------------------------------------
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shp = $slide.Shapes.Item($i)
        if ($shp.HasTable) {
            $shp.Table.ApplyStyle("{8F568475-D16F-4FD6-A70D-62D4BDE3F239}")
        }
    }
}

# ---------------------------------------------------------------------------
# 2) Switch the presentation's colour theme from "Integral" (Red Violet) back
#    to the standard "Office Theme" (Office) colour scheme - i.e. apply a new
#    Design/Theme colour set to the deck.
# ---------------------------------------------------------------------------
$officeColors = @{
    1  = 0         # dk1      000000
    2  = 16777215  # lt1      FFFFFF
    3  = 6968388   # dk2      44546A
    4  = 15132391  # lt2      E7E6E6
    5  = 13998939  # accent1  5B9BD5
    6  = 3243501   # accent2  ED7D31
    7  = 10855845  # accent3  A5A5A5
    8  = 49407     # accent4  FFC000
    9  = 12874308  # accent5  4472C4
    10 = 4697456   # accent6  70AD47
    11 = 12673797  # hlink    0563C1
    12 = 7491477   # folHlink 954F72
}

$slide1 = $p.Slides.Item(1)
$scheme = $slide1.ThemeColorScheme
for ($idx = 1; $idx -le $scheme.Count; $idx++) {
    $scheme.Colors($idx).RGB = $officeColors[$idx]
}
